# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Haba" (Femacal de La Calera) as the
# new row 146, pushing the existing rows 146..178 down to 147..179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 146 - this shifts rows
# 146..178 down to 147..179 (dimension grows from R178 to R179).
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new weekly record.
$ws.Cells.Item(146, 1).Value = 3
$ws.Cells.Item(146, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(146, 3).Value = "Coquimbo"
$ws.Cells.Item(146, 4).Value = 44785
$ws.Cells.Item(146, 5).Value = 5
$ws.Cells.Item(146, 6).Value = 100112026
$ws.Cells.Item(146, 7).Value = "Haba"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 82
$ws.Cells.Item(146, 11).Value = 16000
$ws.Cells.Item(146, 12).Value = 17000
$ws.Cells.Item(146, 13).Value = 16488
$ws.Cells.Item(146, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(146, 16).Value = 660
$ws.Cells.Item(146, 17).Value = 25
$ws.Cells.Item(146, 18).Value = "Hortaliza"
